$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full 37x7 metadata table (header + 36 data rows) as a 2D array
$arr = New-Object 'object[,]' 37,7
$arr[0,0] = 'indicator'
$arr[0,1] = 'description'
$arr[0,2] = 'class'
$arr[0,3] = 'type'
$arr[0,4] = 'unit'
$arr[0,5] = 'sources'
$arr[0,6] = 'update_frequency'
$arr[1,0] = 'name_short'
$arr[1,1] = 'Country, area, or territory name'
$arr[1,2] = 'Character'
$arr[1,3] = 'Direct figure'
$arr[1,4] = 'N/A'
$arr[1,5] = 'WHO'
$arr[1,6] = 'Periodically'
$arr[2,0] = 'iso'
$arr[2,1] = 'ISO-3 code'
$arr[2,2] = 'Character'
$arr[2,3] = 'Direct figure'
$arr[2,4] = 'N/A'
$arr[2,5] = 'ISO'
$arr[2,6] = 'Periodically'
$arr[3,0] = 'who_region'
$arr[3,1] = 'WHO region'
$arr[3,2] = 'Character'
$arr[3,3] = 'Direct figure'
$arr[3,4] = 'N/A'
$arr[3,5] = 'WHO'
$arr[3,6] = 'Periodically'
$arr[4,0] = 'income_group'
$arr[4,1] = 'World Bank 2021 income group classification'
$arr[4,2] = 'Character'
$arr[4,3] = 'Direct figure'
$arr[4,4] = 'N/A'
$arr[4,5] = 'World Bank'
$arr[4,6] = 'Periodically'
$arr[5,0] = 'covax_status'
$arr[5,1] = 'COVAX participation modality'
$arr[5,2] = 'Character'
$arr[5,3] = 'Direct figure'
$arr[5,4] = 'N/A'
$arr[5,5] = 'COVAX Facility'
$arr[5,6] = 'Periodically'
$arr[6,0] = 'covdp_status'
$arr[6,1] = 'COVID-19 Vaccine Delivery Partnership Concerted Support Status'
$arr[6,2] = 'Character'
$arr[6,3] = 'Direct figure'
$arr[6,4] = 'N/A'
$arr[6,5] = 'CoVDP'
$arr[6,6] = 'Periodically'
$arr[7,0] = 'pop'
$arr[7,1] = 'Population, total'
$arr[7,2] = 'Numeric'
$arr[7,3] = 'Direct figure'
$arr[7,4] = 'Individuals'
$arr[7,5] = 'UNPOP'
$arr[7,6] = 'Annually'
$arr[8,0] = 'pop_hcw'
$arr[8,1] = 'Population, total healthcare workers'
$arr[8,2] = 'Numeric'
$arr[8,3] = 'Direct figure'
$arr[8,4] = 'Individuals'
$arr[8,5] = 'ILO Stat'
$arr[8,6] = 'Periodically'
$arr[9,0] = 'pop_older'
$arr[9,1] = 'Population, older adults (as defined by country)'
$arr[9,2] = 'Numeric'
$arr[9,3] = 'Direct figure'
$arr[9,4] = 'Individuals'
$arr[9,5] = 'UNPOP'
$arr[9,6] = 'Annually'
$arr[10,0] = 'date'
$arr[10,1] = 'Date, month, corresponding to data'
$arr[10,2] = 'Date'
$arr[10,3] = 'Direct figure'
$arr[10,4] = 'N/A'
$arr[10,5] = 'N/A'
$arr[10,6] = 'N/A'
$arr[11,0] = 'adm_td_add'
$arr[11,1] = 'Total doses administered, cumulative'
$arr[11,2] = 'Numeric'
$arr[11,3] = 'Direct figure'
$arr[11,4] = 'Doses'
$arr[11,5] = 'WHO COVmart'
$arr[11,6] = 'Monthly'
$arr[12,0] = 'adm_td_month'
$arr[12,1] = 'Total doses administered, net during the month indicated'
$arr[12,2] = 'Numeric'
$arr[12,3] = 'Calculation'
$arr[12,4] = 'Doses'
$arr[12,5] = 'WHO COVmart'
$arr[12,6] = 'Monthly'
$arr[13,0] = 'adm_a1d_add'
$arr[13,1] = 'Individuals having received at least one dose, cumulative'
$arr[13,2] = 'Numeric'
$arr[13,3] = 'Direct figure'
$arr[13,4] = 'Individuals'
$arr[13,5] = 'WHO COVmart'
$arr[13,6] = 'Monthly'
$arr[14,0] = 'adm_a1d_month'
$arr[14,1] = 'Individuals having received at least one dose, net during the month indicated'
$arr[14,2] = 'Numeric'
$arr[14,3] = 'Calculation'
$arr[14,4] = 'Individuals'
$arr[14,5] = 'WHO COVmart'
$arr[14,6] = 'Monthly'
$arr[15,0] = 'adm_fv_add'
$arr[15,1] = 'Individuals having received a complete primary series, cumulative'
$arr[15,2] = 'Numeric'
$arr[15,3] = 'Direct figure'
$arr[15,4] = 'Individuals'
$arr[15,5] = 'WHO COVmart'
$arr[15,6] = 'Monthly'
$arr[16,0] = 'adm_fv_month'
$arr[16,1] = 'Individuals having received a complete primary series, net during the month indicated'
$arr[16,2] = 'Numeric'
$arr[16,3] = 'Calculation'
$arr[16,4] = 'Individuals'
$arr[16,5] = 'WHO COVmart'
$arr[16,6] = 'Monthly'
$arr[17,0] = 'adm_booster_add'
$arr[17,1] = 'Individuals having received a first booster, cumulative'
$arr[17,2] = 'Numeric'
$arr[17,3] = 'Direct figure'
$arr[17,4] = 'Individuals'
$arr[17,5] = 'WHO COVmart'
$arr[17,6] = 'Monthly'
$arr[18,0] = 'adm_booster_month'
$arr[18,1] = 'Individuals having received a first booster, net during the month indicated'
$arr[18,2] = 'Numeric'
$arr[18,3] = 'Calculation'
$arr[18,4] = 'Individuals'
$arr[18,5] = 'WHO COVmart'
$arr[18,6] = 'Monthly'
$arr[19,0] = 'cov_total_a1d'
$arr[19,1] = 'Individuals with at least one dose, as percentage of total population'
$arr[19,2] = 'Numeric'
$arr[19,3] = 'Calculation'
$arr[19,4] = 'Percentage'
$arr[19,5] = 'WHO COVmart'
$arr[19,6] = 'Monthly'
$arr[20,0] = 'cov_total_fv'
$arr[20,1] = 'Individuals with a complete primary series, as percentage of total population'
$arr[20,2] = 'Numeric'
$arr[20,3] = 'Calculation'
$arr[20,4] = 'Percentage'
$arr[20,5] = 'WHO COVmart'
$arr[20,6] = 'Monthly'
$arr[21,0] = 'cov_total_booster'
$arr[21,1] = 'Individuals having received a booster / additional dose, as percentage of total population'
$arr[21,2] = 'Numeric'
$arr[21,3] = 'Calculation'
$arr[21,4] = 'Percentage'
$arr[21,5] = 'WHO COVmart'
$arr[21,6] = 'Monthly'
$arr[22,0] = 'adm_hcw_a1d_add'
$arr[22,1] = 'Healthcare workers having received at least one dose, cumulative'
$arr[22,2] = 'Numeric'
$arr[22,3] = 'Calculation'
$arr[22,4] = 'Individuals'
$arr[22,5] = 'WHO WIISEmart'
$arr[22,6] = 'Monthly'
$arr[23,0] = 'adm_hcw_fv_add'
$arr[23,1] = 'Healthcare workers having received a complete primary series, cumulative'
$arr[23,2] = 'Numeric'
$arr[23,3] = 'Calculation'
$arr[23,4] = 'Individuals'
$arr[23,5] = 'WHO WIISEmart'
$arr[23,6] = 'Monthly'
$arr[24,0] = 'adm_hcw_booster_add'
$arr[24,1] = 'Healthcare workers having received a first booster dose, cumulative'
$arr[24,2] = 'Numeric'
$arr[24,3] = 'Calculation'
$arr[24,4] = 'Individuals'
$arr[24,5] = 'WHO WIISEmart'
$arr[24,6] = 'Monthly'
$arr[25,0] = 'cov_hcw_a1d'
$arr[25,1] = 'Healthcare workers with at least one dose, as percentage of total healthcare worker population'
$arr[25,2] = 'Numeric'
$arr[25,3] = 'Calculation'
$arr[25,4] = 'Percentage'
$arr[25,5] = 'WHO WIISEmart'
$arr[25,6] = 'Monthly'
$arr[26,0] = 'cov_hcw_fv'
$arr[26,1] = 'Healthcare workers with a complete primary series (adjusted), as percentage of total healthcare worker population'
$arr[26,2] = 'Numeric'
$arr[26,3] = 'Calculation'
$arr[26,4] = 'Percentage'
$arr[26,5] = 'WHO WIISEmart'
$arr[26,6] = 'Monthly'
$arr[27,0] = 'cov_hcw_booster'
$arr[27,1] = 'Healthcare workers with a first booster dose, as percentage of total healthcare worker population'
$arr[27,2] = 'Numeric'
$arr[27,3] = 'Calculation'
$arr[27,4] = 'Percentage'
$arr[27,5] = 'WHO WIISEmart'
$arr[27,6] = 'Monthly'
$arr[28,0] = 'adm_old_a1d_add'
$arr[28,1] = 'Older adults having received at least one dose, cumulative'
$arr[28,2] = 'Numeric'
$arr[28,3] = 'Calculation'
$arr[28,4] = 'Individuals'
$arr[28,5] = 'WHO WIISEmart'
$arr[28,6] = 'Monthly'
$arr[29,0] = 'adm_old_fv_add'
$arr[29,1] = 'Older adults having received a complete primary series, cumulative'
$arr[29,2] = 'Numeric'
$arr[29,3] = 'Calculation'
$arr[29,4] = 'Individuals'
$arr[29,5] = 'WHO WIISEmart'
$arr[29,6] = 'Monthly'
$arr[30,0] = 'adm_old_booster_add'
$arr[30,1] = 'Older adults having received a first booster dose, cumulative'
$arr[30,2] = 'Numeric'
$arr[30,3] = 'Calculation'
$arr[30,4] = 'Individuals'
$arr[30,5] = 'WHO WIISEmart'
$arr[30,6] = 'Monthly'
$arr[31,0] = 'cov_old_a1d'
$arr[31,1] = 'Older adults with at least one dose, as percentage of older adult population'
$arr[31,2] = 'Numeric'
$arr[31,3] = 'Calculation'
$arr[31,4] = 'Percentage'
$arr[31,5] = 'WHO WIISEmart'
$arr[31,6] = 'Monthly'
$arr[32,0] = 'cov_old_fv'
$arr[32,1] = 'Older adults with a complete primary series, as percentage of older adult population'
$arr[32,2] = 'Numeric'
$arr[32,3] = 'Calculation'
$arr[32,4] = 'Percentage'
$arr[32,5] = 'WHO WIISEmart'
$arr[32,6] = 'Monthly'
$arr[33,0] = 'cov_old_booster'
$arr[33,1] = 'Older adults with a first booster dose, as percentage of older adult population'
$arr[33,2] = 'Numeric'
$arr[33,3] = 'Calculation'
$arr[33,4] = 'Percentage'
$arr[33,5] = 'WHO WIISEmart'
$arr[33,6] = 'Monthly'
$arr[34,0] = 'del_dose_add'
$arr[34,1] = 'Vaccine doses received, cumulative'
$arr[34,2] = 'Numeric'
$arr[34,3] = 'Direct figure'
$arr[34,4] = 'Doses'
$arr[34,5] = 'UNICEF MDB'
$arr[34,6] = 'Monthly'
$arr[35,0] = 'del_dose_month'
$arr[35,1] = 'Vaccine doses received, net during the month indicated'
$arr[35,2] = 'Numeric'
$arr[35,3] = 'Calculation'
$arr[35,4] = 'Doses'
$arr[35,5] = 'UNICEF MDB'
$arr[35,6] = 'Monthly'
$arr[36,0] = 'est_stock'
$arr[36,1] = 'Estimated remaining vaccine supply, calculated by subtracting adm_td_add from del_dose_add on a monthly basis'
$arr[36,2] = 'Numeric'
$arr[36,3] = 'Calculation'
$arr[36,4] = 'Doses'
$arr[36,5] = 'WHO COVmart; UNICEF MDB'
$arr[36,6] = 'Monthly'

$ws.Range("A1:G37").Value = $arr

# Resize the native table (ListObject) to cover the new data range
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G37"))

# Re-apply column widths (closest achievable values to the authored widths)
$ws.Columns.Item(1).ColumnWidth = 21.166666666666668
$ws.Columns.Item(2).ColumnWidth = 103.66666666666667
$ws.Columns.Item(3).ColumnWidth = 8.0
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 10.166666666666666
$ws.Columns.Item(6).ColumnWidth = 24.5
$ws.Columns.Item(7).ColumnWidth = 17.666666666666668

# Reset selection back to the top-left cell
$ws.Range("A1").Select()

Write-Host "Applied metadata_timeseries_country_month update"
